# Applies the "Added run for config -Min1-Max<=2 all countries" edit:
#   1. Update the section header in A25 to the new, more specific label.
#   2. Append two new result rows (34 & 35) with the Min1_Top1_AllCountry /
#      Min3_Top1_AllCountry metrics.
#   3. Move the active selection to B23 (matches the author's last selection
#      before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the "Only top_n countries..." section header.
$ws.Range("A25").Value = "Only top_n countries (and tied) alowed - First Five Countries"

# 2) Add the two new data rows under the table (row 34 = Min1_Top1_AllCountry,
#    row 35 = Min3_Top1_AllCountry). Columns C:AF hold Recall/Precision/F2
#    triples for each of the 10 categories; AG is left blank (present, empty).

$ws.Range("B34").Value = "Min1_Top1_AllCountry"
$row34 = New-Object 'object[,]' 1,30
$row34[0,0]  = 0.47560975609756101
$row34[0,1]  = 0.185714285714285
$row34[0,2]  = 0.78
$row34[0,3]  = 0.49757281553397997
$row34[0,4]  = 0.19339622641509399
$row34[0,5]  = 0.82
$row34[0,6]  = 0.51094890510948898
$row34[0,7]  = 0.199052132701421
$row34[0,8]  = 0.84
$row34[0,9]  = 0.47706422018348599
$row34[0,10] = 0.168284789644012
$row34[0,11] = 0.88135593220338904
$row34[0,12] = 0.418454935622317
$row34[0,13] = 0.14233576642335699
$row34[0,14] = 0.8125
$row34[0,15] = 0.42857142857142799
$row34[0,16] = 0.15207373271889399
$row34[0,17] = 0.78571428571428503
$row34[0,18] = 0.53738317757009302
$row34[0,19] = 0.20535714285714199
$row34[0,20] = 0.90196078431372495
$row34[0,21] = 0.51229508196721296
$row34[0,22] = 0.19841269841269801
$row34[0,23] = 0.84745762711864303
$row34[0,24] = 0.49118387909319899
$row34[0,25] = 0.194029850746268
$row34[0,26] = 0.79591836734693799
$row34[0,27] = 0.48625792811839302
$row34[0,28] = 0.176245210727969
$row34[0,29] = 0.86792452830188604
$ws.Range("C34:AF34").Value = $row34

$ws.Range("B35").Value = "Min3_Top1_AllCountry"
$row35 = New-Object 'object[,]' 1,30
$row35[0,0]  = 0.29255319148936099
$row35[0,1]  = 0.16176470588235201
$row35[0,2]  = 0.36666666666666597
$row35[0,3]  = 0.25280898876404401
$row35[0,4]  = 0.13636363636363599
$row35[0,5]  = 0.32142857142857101
$row35[0,6]  = 0.39408866995073799
$row35[0,7]  = 0.22535211267605601
$row35[0,8]  = 0.48484848484848397
$row35[0,9]  = 0.39432176656151402
$row35[0,10] = 0.15151515151515099
$row35[0,11] = 0.65789473684210498
$row35[0,12] = 0.36170212765957399
$row35[0,13] = 0.153153153153153
$row35[0,14] = 0.54838709677419295
$row35[0,15] = 0.34574468085106302
$row35[0,16] = 0.180555555555555
$row35[0,17] = 0.44827586206896503
$row35[0,18] = 0.34403669724770602
$row35[0,19] = 0.159574468085106
$row35[0,20] = 0.483870967741935
$row35[0,21] = 0.30232558139534799
$row35[0,22] = 0.14285714285714199
$row35[0,23] = 0.41935483870967699
$row35[0,24] = 0.28985507246376802
$row35[0,25] = 0.13793103448275801
$row35[0,26] = 0.4
$row35[0,27] = 0.28985507246376802
$row35[0,28] = 0.13186813186813101
$row35[0,29] = 0.41379310344827502
$ws.Range("C35:AF35").Value = $row35

# The source columns C,F,I,L,O,R,U,X,AA,AB,AD,AE,AG carry a column-level
# number-format style; the authored rows were written as plain (unstyled)
# values, so strip any inherited formatting back to the default style.
$ws.Range("C34:AG35").Style = "Normal"

# AG34/AG35 stay present-but-empty in the source sheet (trailing blank cell
# in the row's used range).
$ws.Range("AG34").ClearContents()
$ws.Range("AG35").ClearContents()
$ws.Range("AG34:AG35").Style = "Normal"

# 3) Restore the author's final selection/cursor position.
$ws.Range("B23").Select()
